$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title/timestamp cell (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 18:52"

# --- Update country data table ---
# Row 4
$ws.Range("B4").Value = 508126
$ws.Range("C4").Value = 5250
$ws.Range("E4").Value = 460075
$ws.Range("G4").Value = 1080
$ws.Range("H4").Value = 19827
# Row 8
$ws.Range("B8").Value = 123826
$ws.Range("C8").Value = 1655
$ws.Range("E8").Value = 67177
# Row 26
$ws.Range("B26").Value = 7257
$ws.Range("C26").Value = 96
$ws.Range("D26").Value = 411
$ws.Range("E26").Value = 6531
$ws.Range("F26").Value = 184
$ws.Range("G26").Value = 18
$ws.Range("H26").Value = 315
# Row 28
$ws.Range("B28").Value = 6403
$ws.Range("C28").Value = 89
$ws.Range("E28").Value = 6254
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 117
# Row 44
$ws.Range("B44").Value = 3270
$ws.Range("C44").Value = 47
$ws.Range("E44").Value = 2708
$ws.Range("G44").Value = 8
$ws.Range("H44").Value = 62
# Row 112
$ws.Range("D112").Value = 60
$ws.Range("E112").Value = 179
# Row 113
$ws.Range("A113").Value = "Isla de Man"
$ws.Range("B113").Value = 226
$ws.Range("C113").Value = 25
$ws.Range("D113").Value = 112
$ws.Range("E113").Value = 113
$ws.Range("F113").Value = 11
$ws.Range("H113").Value = 1
# Row 114
$ws.Range("A114").Value = "Consejo Danes para los Refugiados"
$ws.Range("B114").Value = 223
$ws.Range("D114").Value = 16
$ws.Range("E114").Value = 187
$ws.Range("H114").Value = 20
# Row 115
$ws.Range("A115").Value = "Guinea"
$ws.Range("B115").Value = 212
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 15
$ws.Range("E115").Value = 197
$ws.Range("F115").Value = 0
$ws.Range("H115").Value = 0
# Row 127
$ws.Range("B127").Value = 129
$ws.Range("C127").Value = 2
$ws.Range("D127").Value = 84
$ws.Range("E127").Value = 45
# Row 168
$ws.Range("A168").Value = "Sudan"
$ws.Range("C168").Value = 2
# Row 169
$ws.Range("A169").Value = "Angola"
$ws.Range("C169").Value = 0
# Row 183
$ws.Range("A183").Value = "Belice"
$ws.Range("C183").Value = 3
$ws.Range("E183").Value = 11
$ws.Range("F183").Value = 1
$ws.Range("H183").Value = 2
# Row 184
$ws.Range("A184").Value = "Zimbabue"
$ws.Range("B184").Value = 13
$ws.Range("E184").Value = 10
$ws.Range("H184").Value = 3
# Row 185
$ws.Range("A185").Value = "San Cristobal y Nieves"
$ws.Range("D185").Value = 0
$ws.Range("E185").Value = 12
# Row 186
$ws.Range("A186").Value = "San Vicente y las Granadinas"
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 1
$ws.Range("E186").Value = 11
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0
# Row 187
$ws.Range("A187").Value = "Malaui"
$ws.Range("C187").Value = 3
$ws.Range("D187").Value = 0
$ws.Range("E187").Value = 10
$ws.Range("F187").Value = 1
$ws.Range("G187").Value = 1
$ws.Range("H187").Value = 2
# Row 188
$ws.Range("A188").Value = "Suazilandia"
$ws.Range("B188").Value = 12
$ws.Range("D188").Value = 7
$ws.Range("E188").Value = 5
# Row 189
$ws.Range("A189").Value = "Seychelles"
$ws.Range("D189").Value = 0
$ws.Range("E189").Value = 11
# Row 190
$ws.Range("A190").Value = "Republica del Chad"
$ws.Range("D190").Value = 2
$ws.Range("E190").Value = 9
# Row 191
$ws.Range("A191").Value = "Groenlandia"
$ws.Range("B191").Value = 11
$ws.Range("D191").Value = 11
$ws.Range("E191").Value = 0
$ws.Range("F191").Value = 0
$ws.Range("H191").Value = 0
# Row 197
$ws.Range("A197").Value = "Nicaragua"
$ws.Range("C197").Value = 1
# Row 198
$ws.Range("A198").Value = "Islas Turcas y Caicos"
$ws.Range("C198").Value = 0
# Row 209
$ws.Range("A209").Value = "Burundi"
# Row 210
$ws.Range("A210").Value = "Anguila"
# Row 212
$ws.Range("A212").Value = "Papua Nueva Guinea"
# Row 213
$ws.Range("A213").Value = "Bonaire, San Eustaquio y Saba"
